# "fixed export and fixing maps"
#
# The Borjomi Municipality "Area" sheet gets simplified:
#   - the "(according to the population census data)" subtitle row is removed
#   - the 1989 / 2002 historical columns are removed, leaving only the 2014
#     column (which slides left into column B)
#   - remaining rows are given a slightly taller custom row height
#
# Row/column layout before -> after:
#   row1 "Area of the municipality of Borjomi"      -> row1 (unchanged)
#   row2 "(according to the population census data)" -> REMOVED
#   row3 (blank spacer)                              -> row2
#   row4 "(sq. km)"                                   -> row3
#   row5 blank | 1989 | 2002 | 2014                   -> row4 blank | 2014
#   row6 "Area" | 1189 | 1189 | 1189                  -> row5 "Area" | 1189

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the census-data subtitle row entirely (old row 2).
$ws.Rows.Item(2).Delete() | Out-Null

# Remove the now-unwanted 1989 and 2002 columns (old columns B and C);
# the 2014 column shifts left into column B, carrying its own values and
# formatting with it.
$ws.Columns.Item(2).Delete() | Out-Null
$ws.Columns.Item(2).Delete() | Out-Null

# The remaining five rows get a slightly taller, explicit row height.
$ws.Rows("1:5").RowHeight = 20.1

# Match the saved cursor/selection position recorded in the workbook.
$ws.Range("H19").Select() | Out-Null

$wb.Save()
